$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("M2").Value = 3.832616
$ws.Range("N2").Value = 11.497848
$ws.Range("O2").Value = 0.264217765931355
$ws.Range("P2").Value = 0.264217765931355
$ws.Range("Q2").Value = 0.2606817649333333
$ws.Range("R2").Value = 2.3461358844
$ws.Range("S2").Value = 0.264217765931355
$ws.Range("T2").Value = 0.264217765931355

# Row 3 updates
$ws.Range("O3").Value = 0.6031799756961482
$ws.Range("P3").Value = 0.6031799756961482
$ws.Range("S3").Value = 0.6031799756961482
$ws.Range("T3").Value = 0.6031799756961482

# Row 4 updates
$ws.Range("M4").Value = 1.255882
$ws.Range("N4").Value = 3.767646
$ws.Range("O4").Value = 0.08657959375878042
$ws.Range("P4").Value = 0.08657959375878042
$ws.Range("Q4").Value = 0.08542090736666667
$ws.Range("R4").Value = 0.7687881663
$ws.Range("S4").Value = 0.08657959375878042
$ws.Range("T4").Value = 0.08657959375878042

# Row 5 updates
$ws.Range("M5").Value = 0.6675826666666667
$ws.Range("N5").Value = 2.002748
$ws.Range("O5").Value = 0.04602266461371635
$ws.Range("P5").Value = 0.04602266461371635
$ws.Range("Q5").Value = 0.04540674771111111
$ws.Range("R5").Value = 0.4086607294
$ws.Range("S5").Value = 0.04602266461371635
$ws.Range("T5").Value = 0.04602266461371635
